$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-6 (columns A-J). A new row (row 6) is inserted and
# the values for rows 2-5 are shifted/updated per the updated training
# schedule.
$data = @(
    @(1, 9, 7, 6, 4, -3, -3, 43, 5, "train_dim2_1"),
    @(2, 5, 5, 0, 4, -5, -1, 65, 5, "train_dim2_1"),
    @(3, 8, 6, 7, 1, -1, -5, 21, 5, "train_dim2_1"),
    @(4, 5, 7, 1, 5, -4, -2, 54, 5, "train_dim2_1"),
    @(5, 9, 6, 7, 2, -2, -4, 32, 5, "train_dim2_1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
}

# Update the selection to match the recorded state (I1 selected).
$ws.Range("I1").Select()
